$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Remove the "Lương" sheet (payroll report no longer generated here)
$wb.Worksheets.Item('Lương').Delete()

# 2. Update "Đơn 1 bác sĩ" sheet: add "Nhóm dịch vụ" column + 4 discount columns,
#    and add a new data row, recomputing the totals row.
$ws = $wb.Worksheets.Item('Đơn 1 bác sĩ')

# Header row (row 1)
$ws.Cells.Item(1, 1).Value = 'Tiền tố'
$ws.Cells.Item(1, 2).Value = 'Mã dịch vụ'
$ws.Cells.Item(1, 3).Value = 'Ngày thực hiện'
$ws.Cells.Item(1, 4).Value = 'Cơ sở'
$ws.Cells.Item(1, 5).Value = 'Khách hàng'
$ws.Cells.Item(1, 6).Value = 'Nguồn khách'
$ws.Cells.Item(1, 7).Value = 'Nhóm dịch vụ'
$ws.Cells.Item(1, 8).Value = 'Tên dịch vụ'
$ws.Cells.Item(1, 9).Value = 'Sale chính'
$ws.Cells.Item(1, 10).Value = 'Đơn giá gốc'
$ws.Cells.Item(1, 11).Value = 'Sale phụ'
$ws.Cells.Item(1, 12).Value = 'Upsale'
$ws.Cells.Item(1, 13).Value = 'Đơn giá'
$ws.Cells.Item(1, 14).Value = 'Thanh toán lần đầu'
$ws.Cells.Item(1, 15).Value = 'Trả sau'
$ws.Cells.Item(1, 16).Value = 'Đã thanh toán'
$ws.Cells.Item(1, 17).Value = 'Dư nợ'
$ws.Cells.Item(1, 18).Value = 'Bác sĩ 1'
$ws.Cells.Item(1, 19).Value = 'Bác sĩ 2'
$ws.Cells.Item(1, 20).Value = 'Phụ phẫu 1'
$ws.Cells.Item(1, 21).Value = 'Phụ phẫu 2'
$ws.Cells.Item(1, 22).Value = 'Công phụ phẫu 1'
$ws.Cells.Item(1, 23).Value = 'Công phụ phẫu 2'
$ws.Cells.Item(1, 24).Value = 'Tỉ lệ chiết khấu sale chính'
$ws.Cells.Item(1, 25).Value = 'Tỉ lệ chiết khấu sale phụ'
$ws.Cells.Item(1, 26).Value = 'Chiết khấu sale chính'
$ws.Cells.Item(1, 27).Value = 'Chiết khấu sale phụ'

# Row 2 (existing order, re-laid-out across the new columns)
$ws.Cells.Item(2, 1).Value = 'HD-LUXURY'
$ws.Cells.Item(2, 2).Value = 521
$ws.Cells.Item(2, 3).Value = "'07-05-2024"
$ws.Cells.Item(2, 4).Value = 'SÓC TRĂNG'
$ws.Cells.Item(2, 5).Value = 'đa ni '
$ws.Cells.Item(2, 6).Value = 'CTV'
$ws.Cells.Item(2, 7).Value = 'Môi'
$ws.Cells.Item(2, 8).Value = 'Phun môi'
$ws.Cells.Item(2, 9).Value = 'Thạch Hoàng Nhân'
$ws.Cells.Item(2, 10).Value = 5500000
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 5500000
$ws.Cells.Item(2, 14).Value = 5500000
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 5500000
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 'Bác Sĩ Ngoài'
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 20).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 0
$ws.Cells.Item(2, 25).Value = 0
$ws.Cells.Item(2, 26).Value = 0
$ws.Cells.Item(2, 27).Value = 0

# Row 3 (new order)
$ws.Cells.Item(3, 1).Value = 'HD-LUXURY'
$ws.Cells.Item(3, 2).Value = 537
$ws.Cells.Item(3, 3).Value = "'07-11-2024"
$ws.Cells.Item(3, 4).Value = 'SÓC TRĂNG'
$ws.Cells.Item(3, 5).Value = 'dương thị lệ '
$ws.Cells.Item(3, 6).Value = 'Khách cũ'
$ws.Cells.Item(3, 7).Value = 'Vùng mắt'
$ws.Cells.Item(3, 8).Value = 'Phun mày'
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 1000000
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 1000000
$ws.Cells.Item(3, 14).Value = 1000000
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 1000000
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 'Bác Sĩ Ngoài'
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 20).Value = 0
$ws.Cells.Item(3, 21).Value = 0
$ws.Cells.Item(3, 22).Value = 0
$ws.Cells.Item(3, 23).Value = 0
$ws.Cells.Item(3, 24).Value = 0.1
$ws.Cells.Item(3, 25).Value = 0
$ws.Cells.Item(3, 26).Value = 100000
$ws.Cells.Item(3, 27).Value = 0

# Row 4 ("Tổng" totals row)
$ws.Cells.Item(4, 1).Value = 'Tổng'
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = ''
$ws.Cells.Item(4, 4).Value = ''
$ws.Cells.Item(4, 5).Value = ''
$ws.Cells.Item(4, 6).Value = ''
$ws.Cells.Item(4, 7).Value = ''
$ws.Cells.Item(4, 8).Value = ''
$ws.Cells.Item(4, 9).Value = ''
$ws.Cells.Item(4, 10).Value = 6500000
$ws.Cells.Item(4, 11).Value = ''
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 6500000
$ws.Cells.Item(4, 14).Value = 6500000
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 6500000
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = ''
$ws.Cells.Item(4, 19).Value = ''
$ws.Cells.Item(4, 20).Value = ''
$ws.Cells.Item(4, 21).Value = ''
$ws.Cells.Item(4, 22).Value = 0
$ws.Cells.Item(4, 23).Value = 0
$ws.Cells.Item(4, 24).Value = 0.1
$ws.Cells.Item(4, 25).Value = 0
$ws.Cells.Item(4, 26).Value = 100000
$ws.Cells.Item(4, 27).Value = 0

Write-Host "Done"
